$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "Introduction"
$ws.Range("C1").Value = "Show head"
$ws.Range("E1").Value = "High risk patients"
$ws.Range("G1").Value = "Cross tab"
$ws.Range("I1").Value = "Scatter plot"

# --- Row 2 ---
$ws.Range("A2").Value = "Hi there, this is a tool I created to demonstrate my program to compute Framingham risk scores and demonstrate ways to analyze them. These scores estimate the person's risk of coronary heart disease in the coming decade. They are based on the participant's age, sex, systolic blood pressure, if they are being treated for high blood pressure, smoking status, HDL cholesterol and total cholesterol. Because the weighting schemes used in the calculation can be cumbersome since they vary based on factors like sex and age, so it's useful to automate the context so the scores can be applied to a large dataset.`r`n`r`nThis tool has several demonstrations to run through. Enter their number to run through."
$ws.Range("C2").Value = "Here's what the data looks like at the start. All it has is some basic hypothetical health record data with 250 cases. This is everything needed to calculate Framingham scores."
$ws.Range("E2").Value = "Let's put these data to work. `r`n`r`nThere are {len(atrisk)} people in the dataset with over 20% risk. `r`n`r`nPerhaps we should review them so we can notify their care team (HIPAA notwithstanding)? `r`n`r`nPress enter to see their IDs."
$ws.Range("G2").Value = "It might be informative to examine if people's forecasted risk is associated with receiving SNAP benefits. Since SNAP qualification hinges (in part) on socioeconomic status it may be a useful indicator in general, but this would help policy advocates to decide if they should spend time pushing for that new proposal to provide low cost cardiac evaluations for SNAP recipients. `r`n`r`nPress enter to show the cross tabulation."
$ws.Range("I2").Value = "Do Framingham risk estimates correlate with an individual's annual medical expenditures? It would make sense if they did- a high risk person might be engaging more with their doctors to begin with, or a person with lots of medical expenses might have more stress and less means/bandwidth to sustain a healthy diet, inflating some parts of the Framingham scores (especially systolic blood Enterure and cholesterol).`r`n`r`nPress enter to start the analysis with some histograms and visuals, then finishing with some regression models."

# --- Row 3 ---
$ws.Range("A3").Value = "****Interact with the tool:`r`n1) Illustrate the steps of the score calculator`r`n2) Enter information to get a risk score estimate back`r`n****Run Analyses on a fake dataset of 250 cases.`r`n3) Identify high risk patients.`r`n4) Examine the association between SNAP receipt and risk scores.`r`n5) Use linear regression to model the relationships between risk scores and medical spending.`r`n6) Close the program"
$ws.Range("C3").Value = "First, patient's information is converted into scores, where higher scores=higher risk. For instance, some people are fairly young, so they have scores of 0, but row 2 is a 55 year old woman, who gets a value of 8 for age.`r`n`r`nThe scoring criteria vary based on several factors. As one example- For total cholesterol, the score values are different for women, and they vary yet again by age."
$ws.Range("G3").Value = "NOTE: These numbers are randomly generated, yours will differ a bit.`r`n`r`nLooking at my output,these results indicate that people who are currently on SNAP are more likely to have high risk scores. While 5.32% of people off of SNAP were high risk, 18.87% of people on SNAP were. That's over three times higher risk, although the high risk cases are still in a substantial minority in both groups. The chi square test is statistically significant in this case, which would typically suggest that the result we are observing is unlikely to be due to random sampling error, but given that our sample is non-random we probably should be careful about using this test to make a more general claim. The Cramer's V statistic of .19 is typically interpreted as evidence of only a weak association. It's still possible this result is just random luck of the draw and that the correlation is fairly small, but ultimately it is still a substantively significant difference, as it demonstrates a  clear actionable disparity. `r`n`r`nA caution- chi-square tests cannot be used with less than 4 cases per cell. Here, there are actually 5 in one cell if you check it out. A fix might be to lower the threshold from 20% to 18% risk, though obviously we lower the apparent urgency of our analysis. This would free up a few extra cases.`r`n`r`nA subsequent step might involve running a logistic regression model to further investigate this association, adjusting for other covariates that might confound the relationship. As just one example, SNAP's link to risk scores might actually be explained by chronic financial stress, which would explain both SNAP receipt AND high cardiac risk (especially due to blood pressure). In general, we would want to at least control for as many demographics as possible, especially ones other than age and sex, since including those two might induce biased estimates/standard errors in light of multicollinearity due to them being used in the computation of Framingham scores. A second strategy would be to try to establish a predictive model via a machine learning algorithm, such as a random forest model. Given that we have only 200 cases, we should probably hold off on the latter approach until we have more data, however. "
$ws.Range("I3").Value = "These histograms visualize which values are most common in our data. A tall bar means more people had responses in that range. Looking at the histogram from Framingham scores, three things stick out:`r`n`r`n1) The skew of the distribution is apparent- most respondents pile into lower risk scores (thankfully), with high values becoming comparatively more rare. The gap between the red line, the mean, and the blue line, the median, further illustrates this. `r`n`r`n2) The strange emptiness in some of the higher categories. It's because of how raw scores are translated into risk estimates- as one example, in men a raw score of 14 for men turns into a risk estimate of  16, but a raw score of 15 goes straight to 20. The skewness of this distribution might give us some concern if this was an outcome variable, given the assumptions made in many analyses about normality of the dependent variable. However, independent variables are more flexible- at risk of getting too technical, what matters in the end is that our residual error is normally distributed, not our independent variables themselves.`r`n`r`n3) The overall high mean for risk scores. Most research I looked at indicates that there should be even fewer people in the higher risk values. But, remember what we saw about SNAP in the other analysis- almost half our sample is SNAP recipients, so their heightned risk is probably reflected here.`r`n`r`nThe distribution for medical expenditures is far less concerning- it follows an ideal normal distribution with a mean of PATIENTMEAN and a standard deviation of PATIENTSTD. 95% cases fall within +/- STDRANGE. There is some piling in a few spots, but it is still fairly symmetrial. In my graph there might also be a few outliers, too, with a few cases with particularly low scores. These should be checked up on to make sure they are valid entries, and stakeholders might need to be spoken to to help with contextualizing the results.`r`n`r`nPress enter to continue the analysis.`r`n"
$ws.Range("B3").Value = ""
$ws.Range("B3").WrapText = $true

# --- Row 4 ---
$ws.Range("A4").Value = "`r`nNOTE: The scoring system used are from:`r`n Sullivan LM, Massaro JM, D'Agostino RB Sr. Presentation of multivariate data for clinical use: The Framingham Study risk score functions. Stat Med. 2004 May 30;23(10):1631-60. doi: 10.1002/sim.1742. PMID: 15122742.`r`nAll data are FAKE and correlations are manufactured. Do NOT draw any real conclusions from these analyses."
$ws.Range("C4").Value = "Raw Framingham can (and often should) be converted into their forecasted risk of cardiovascular disease in the next decade. If we were to depend on just raw scores we would deceive ourselves- the calculation,by design, gives women higher scores on average, BUT those values correspond to lower risk. To illustrate: A woman with a score of 20 has an 11% 10-year risk. But for a man with a score of 20 has a risk greater than 30%!  Be aware, though, that some have criticized Framingham scores for underestimating women's risk of heart disease, in part due to changes in our undestanding of how heart diseases present themselves in female bodies, which were previously assumed to be the same as male bodies.**`r`n`r`nHere, you can see the raw Framingham scores, sex, and risk level for a few cases side by side. I also created a few handy categorical flags.`r`n`r`n** `r`nPark KE, Pepine CJ. Assessing cardiovascular risk in women: looking beyond traditional risk factors. Trends Cardiovasc Med. 2015 Feb;25(2):152-3. doi: 10.1016/j.tcm.2014.10.024. Epub 2014 Nov 4. PMID: 25468290; PMCID: PMC4664450.`r`n`r`nDusenbery, M. Doing Harm: The Truth About How Bad Medicine and Lazy Science Leave Women Dismissed, Misdiagnosed, and Sick. HarperOne. 2017."
$ws.Range("I4").Value = "A simple scatter plot helps to clarify the relationship between the two variables. I'll start with a simple one and then we can use the same graph to examine covariates. At the end you can request my interpretation.`r`n`r`nPress enter to display the first scatter plot."

# --- Column I, rows 5-12 ---
$ws.Range("I5").Value = "Here it is- as you can see, it's a pretty loose correlation.`r`n`r`nPress enter for the next graph, which will disaggregate these results by sex."
$ws.Range("I6").Value = "Here's the same association broken down by sex. The orange dots are females, the blue ones are males.`r`n`r`nIf you would like an interpretation of the graphs, enter Y, otherwise enter N."
$ws.Range("I7").Value = "Both of these graphs indicate a positive association between the two variables- as Framingham risk estimates increase, the annual medical spending tends to go up. Though subtle, it is possible that there is a curvilinear trend as opposed to a linear trend, too, which could always be tested for in a later analysis. The pearson's R here is CORRCOEF, indicating a moderate/weak association, but certainly not ignorable. `r`n`r`nThe second graph's breakdown by sex complicates the picture by revealing that, on average, females clearly have higher medical spending than males. Given that earlier we talked about who men tend to have higher risk scores, this is a clear confounder that should be accounted for statistically."
$ws.Range("I8").Value = "In light of those gender differences that we observed, we shouldn't depend on that intial correlation. Instead, we should use a linear regression model to adjust for gender and SNAP membership at once to better isolate the relationship between risk scores on medical expenditures.`r`n`r`nPress enter to run the model."
$ws.Range("I9").Value = "Enter 1 for the same model with standardized coefficients `r`nEnter 2 for a review of the regression diagnostics`r`nEnter 3 for a full interpretation of the model and diagnostics`r`nEnter 4 to move on.`r`n`r`nAll options except for #4 will return you to these options so you can see everything. "
$ws.Range("I10").Value = "MODEL FIT:`r`nThe model's adjusted R-squared value is .389, meaning that these three variables can collectively explain 38.7% of the variance in annual medical spending. Though certainly not guaranteed accuracy, this is a really high number for social sciences data and a model with only two variables!`r`n`r`nCOEFFICIENT ESTIMATES:`r`nThe intercept reflects all covariates=0, reflecting the estimated annual spending for a female with 0 estimated risk. In this case, that's still `$2,327, still rather high!`r`n`r`nRisk_pct's estimate is statistically significant (dubious use for inference in this case), indicating that the result probably shouldn't be chocked up to random luck of the draw. The coefficient tells us that as a person's risk of cardiovascular disease increases by 1%, we expect that their annual medical spending will increase by `$54. Though that sounds small, recall that the risk scores go up to 30, so that's as much as a 54*30=`$1,620 increase. Be mindful, though, that this is a description of our current dataset's patterns, not necessarily a prediction of the future.`r`n`r`nThe coefficient for gender is also statistically significant, estimating that  males in these data spend about `$748 less than females on medical costs per year. `r`n`r`nThe coefficient for SNAP is small and is not statistically significant- it should proably be ignored. Notice the size of the standard error, which is almost equal to the size of the coefficient estimate- we can't conclude if SNAP membership is associated with an increase of decrease in data.  `r`n`r`nIf we were to look at the standardized coefficients, we can confirm that the correlation for risk scores is the strongest in the model, with a 1 std deviation increase in risk % corresponding with a .59 standard deviation increase in medical spending. Contrast that with the -.487 standard deviation change when comparing males to females.`r`n`r`nDIAGNOSTICS:`r`n`r`nThe residual plots give some pause- notice how in some spots the residuals are tightly clustered around the reference line, which in other spots they are more spread out. This is an indication of heteroskedasticity, though the Breusch Pagan test is not significant, suggesting heteroskedasticity is not a massive problem. But the CCPR plot again suggests a shape that is not entirely linear, with a bit of a curve visible. This, and the histogram of errors that is clearly not normal, suggests that this model is still not a great fit to the data. We could try to fix this by introducing more covariates, if possible, or trying different functional forms for the variables we've already got, such as a quadratic term for the Framingham risk pct.`r`n`r`nThe outlier plot indicates six cases that are likely to be outliers. We should run a later model without those cases to see if they are influencing our conclusions at all. "
$ws.Range("I11").Value = "The square term considerably improves the models fit a bit, as evinced by the increased r-squared. But it's also much harder to interpret what is going on now because the relationship between the risk scores and annual medical spending is now curvilinear. We need to visualize the estimates to help us see what is going on.`r`n`r`nPress enter to graph the results."
$ws.Range("I12").Value = "Looking at this graph, we can see that spending is expected to increase sharply for about the first 14% in the risk scale, at which point annual spending tapers out. We can use this regression equation to make particular estimates of a person's expenditures based on their characteristics if we wanted to- just remember to factor gender in, too!`r`n`r`nThat concludes the regression analysis- this is probably the best fitting model we can get with such limited data and a basic linear regression. However, we should be cautious about trusting this model- the omission of many key variables is likely to have biased out models in ways that we cannot know. `r`n`r`n******Press enter to go back to the main hub.******"

# --- Wrap text for the main content cells (rows 2-4 columns A,C,E,G,I and column I rows 5-12) ---
$ws.Range("A2:A4").WrapText = $true
$ws.Range("C3:C4").WrapText = $true
$ws.Range("E2").WrapText = $true
$ws.Range("G2:G3").WrapText = $true
$ws.Range("I2:I12").WrapText = $true

# --- Update the selection / view to match the saved view (E2 selected, scrolled near top) ---
$ws.Range("E2").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 2
} catch {}
